$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 7400.3
$ws.Range("I70").Value = 2000
$ws.Range("J70").Value = 8750.375
$ws.Range("K70").Value = 6000
$ws.Range("L70").Value = 26251.125
$ws.Range("M70").Value = -5730
$ws.Range("N70").Value = -26791.125
$ws.Range("H73").Value = 7400.3
$ws.Range("I73").Value = 2000
$ws.Range("J73").Value = 8750.375
$ws.Range("K73").Value = 6000
$ws.Range("L73").Value = 26251.125
$ws.Range("M73").Value = -5064
$ws.Range("N73").Value = -28123.125
$ws.Range("H88").Value = 2872.5
$ws.Range("I88").Value = 2663.3333
$ws.Range("K88").Value = 2663.3333
$ws.Range("M88").Value = -2257.3333
$ws.Range("H91").Value = 2872.5
$ws.Range("I91").Value = 2663.3333
$ws.Range("K91").Value = 2663.3333
$ws.Range("M91").Value = -1259.3333
$ws.Range("H116").Value = 4579.7144
$ws.Range("I116").Value = 3639.5
$ws.Range("J116").Value = 5833.3335
$ws.Range("K116").Value = 3639.5
$ws.Range("L116").Value = 5833.3335
$ws.Range("M116").Value = -197.5
$ws.Range("N116").Value = -12717.3335
$ws.Range("H125").Value = 2361
$ws.Range("J125").Value = 2361
$ws.Range("L125").Value = 21249
$ws.Range("N125").Value = -26169
$ws.Range("H137").Value = 2158.6191
$ws.Range("I137").Value = 1295.7273
$ws.Range("J137").Value = 3107.8
$ws.Range("K137").Value = 3887.1819
$ws.Range("L137").Value = 9323.400000000001
$ws.Range("M137").Value = -1337.1819
$ws.Range("N137").Value = -14423.4
$ws.Range("H141").Value = 2744.111
$ws.Range("J141").Value = 2516
$ws.Range("L141").Value = 7548
$ws.Range("N141").Value = -17908

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 250.5
$ws.Range("I5").Value = 250.5
$ws.Range("K5").Value = 250.5
$ws.Range("M5").Value = -138.5
$ws.Range("H61").Value = 3082.8333
$ws.Range("I61").Value = 3082.8333
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 3082.8333
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -2870.8333
$ws.Range("N61").ClearContents()
$ws.Range("H88").Value = 2181.25
$ws.Range("I88").Value = 2650
$ws.Range("J88").Value = 1400
$ws.Range("K88").Value = 2650
$ws.Range("L88").Value = 1400
$ws.Range("M88").Value = -2244
$ws.Range("N88").Value = -2212
$ws.Range("H91").Value = 2181.25
$ws.Range("I91").Value = 2650
$ws.Range("J91").Value = 1400
$ws.Range("K91").Value = 2650
$ws.Range("L91").Value = 1400
$ws.Range("M91").Value = -1246
$ws.Range("N91").Value = -4208
$ws.Range("H122").Value = 1699.5
$ws.Range("I122").Value = 1266.3334
$ws.Range("K122").Value = 3799.0002
$ws.Range("M122").Value = -1349.0002
$ws.Range("H132").Value = 1546
$ws.Range("I132").Value = 1584.4445
$ws.Range("J132").Value = 1200
$ws.Range("K132").Value = 4753.333500000001
$ws.Range("L132").Value = 3600
$ws.Range("M132").Value = -2223.333500000001
$ws.Range("N132").Value = -8660
$ws.Range("H133").Value = 148333.33
$ws.Range("J133").Value = 148333.33
$ws.Range("L133").Value = 148333.33
$ws.Range("N133").Value = -153393.33
$ws.Range("H136").Value = 3082.8333
$ws.Range("I136").Value = 3082.8333
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 9248.499899999999
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -6698.499899999999
$ws.Range("N136").ClearContents()

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 250.5
$ws.Range("I4").Value = 250.5
$ws.Range("K4").Value = 250.5
$ws.Range("M4").Value = -135.5
$ws.Range("H20").Value = 2716.4443
$ws.Range("J20").Value = 398.5
$ws.Range("L20").Value = 398.5
$ws.Range("N20").Value = -892.5
$ws.Range("H82").Value = 29139
$ws.Range("I82").Value = 16283.546
$ws.Range("J82").Value = 99844
$ws.Range("K82").Value = 16283.546
$ws.Range("L82").Value = 99844
$ws.Range("M82").Value = -15900.546
$ws.Range("N82").Value = -100610
$ws.Range("H85").Value = 29139
$ws.Range("I85").Value = 16283.546
$ws.Range("J85").Value = 99844
$ws.Range("K85").Value = 16283.546
$ws.Range("L85").Value = 99844
$ws.Range("M85").Value = -14957.546
$ws.Range("N85").Value = -102496
$ws.Range("H134").Value = 3571.261
$ws.Range("I134").Value = 852
$ws.Range("K134").Value = 2556
$ws.Range("M134").Value = -21
$ws.Range("H140").Value = 49184
$ws.Range("J140").Value = 49184
$ws.Range("L140").Value = 49184
$ws.Range("N140").Value = -59544

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 10000000
$ws.Range("J3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("N3").ClearContents()
$ws.Range("H7").Value = 72.09524
$ws.Range("I7").Value = 55
$ws.Range("J7").Value = 87.63636
$ws.Range("K7").Value = 55
$ws.Range("L7").Value = 87.63636
$ws.Range("M7").Value = 58
$ws.Range("N7").Value = -313.63636
$ws.Range("H70").Value = 32250
$ws.Range("J70").Value = 32250
$ws.Range("L70").Value = 32250
$ws.Range("N70").Value = -32880
$ws.Range("H73").Value = 32250
$ws.Range("J73").Value = 32250
$ws.Range("L73").Value = 32250
$ws.Range("N73").Value = -34434
$ws.Range("H122").Value = 1962.5
$ws.Range("J122").Value = 1850
$ws.Range("L122").Value = 5550
$ws.Range("N122").Value = -10450
$ws.Range("H132").Value = 3985.25
$ws.Range("I132").Value = 4059
$ws.Range("J132").Value = 3469
$ws.Range("K132").Value = 12177
$ws.Range("L132").Value = 10407
$ws.Range("M132").Value = -9647
$ws.Range("N132").Value = -15467
$ws.Range("H133").Value = 45326
$ws.Range("J133").Value = 45326
$ws.Range("L133").Value = 45326
$ws.Range("N133").Value = -50386

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H10").Value = 227.875
$ws.Range("J10").Value = 375
$ws.Range("L10").Value = 1125
$ws.Range("N10").Value = -1403
$ws.Range("H13").Value = 244.7
$ws.Range("J13").Value = 280.75
$ws.Range("L13").Value = 842.25
$ws.Range("N13").Value = -1178.25
$ws.Range("H39").Value = 5003.091
$ws.Range("J39").Value = 5174.6665
$ws.Range("L39").Value = 15523.9995
$ws.Range("N39").Value = -16111.9995
$ws.Range("H113").Value = 1305.8889
$ws.Range("J113").Value = 1438.9231
$ws.Range("L113").Value = 4316.7693
$ws.Range("N113").Value = -8656.7693

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 13229.9
$ws.Range("I43").Value = 6679.8
$ws.Range("K43").Value = 6679.8
$ws.Range("M43").Value = -6528.8
$ws.Range("H102").Value = 1701.5625
$ws.Range("I102").Value = 959.5454999999999
$ws.Range("J102").Value = 3334
$ws.Range("K102").Value = 959.5454999999999
$ws.Range("L102").Value = 3334
$ws.Range("M102").Value = 662.4545000000001
$ws.Range("N102").Value = -6578
$ws.Range("H122").Value = 2289.3845
$ws.Range("J122").Value = 3092.8333
$ws.Range("L122").Value = 9278.499899999999
$ws.Range("N122").Value = -14178.4999
$ws.Range("H132").Value = 93291.91
$ws.Range("I132").Value = 102319.7
$ws.Range("K132").Value = 306959.1
$ws.Range("M132").Value = -304429.1
$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1582.15
$ws.Range("I16").Value = 1560.0526
$ws.Range("J16").Value = 2002
$ws.Range("K16").Value = 1560.0526
$ws.Range("L16").Value = 2002
$ws.Range("M16").Value = -1390.0526
$ws.Range("N16").Value = -2342
$ws.Range("H53").Value = 9500
$ws.Range("I53").Value = 9000
$ws.Range("K53").Value = 9000
$ws.Range("M53").Value = -8482
$ws.Range("H55").Value = 1294.9231
$ws.Range("I55").Value = 1121
$ws.Range("K55").Value = 1121
$ws.Range("M55").Value = -948
$ws.Range("H122").Value = 6134.3335
$ws.Range("J122").Value = 4005
$ws.Range("L122").Value = 12015
$ws.Range("N122").Value = -16915
$ws.Range("H123").Value = 86000
$ws.Range("J123").Value = 86000
$ws.Range("L123").Value = 86000
$ws.Range("N123").Value = -95800

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 37885.8
$ws.Range("I41").Value = 33017.332
$ws.Range("K41").Value = 33017.332
$ws.Range("M41").Value = -32627.332
$ws.Range("H122").Value = 2736.7727
$ws.Range("I122").Value = 2563.3684
$ws.Range("K122").Value = 7690.1052
$ws.Range("M122").Value = -5240.1052
$ws.Range("H132").Value = 1614.25
$ws.Range("I132").Value = 1337.1
$ws.Range("K132").Value = 4011.3
$ws.Range("M132").Value = -1481.3
$ws.Range("H141").Value = 230000
$ws.Range("J141").Value = 90000
$ws.Range("L141").Value = 90000
$ws.Range("N141").Value = -100360
